# Sprint 2 standup notes update
#  - Scrum Master name change
#  - Date field change
#  - Status table content refresh for the three team members

$d = $word.ActiveDocument
$wNs = ' xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Scrum Master: Alex Thurston -> Jay Peterson -------------------------
# Scope the Find/Replace to just the "Scrum Master" paragraph so the
# "Alex Thurston" entry still present in the status table (row 1, Team
# Member column) is left untouched.
$scrumMasterPara = $d.Paragraphs.Item(3)
[void]$scrumMasterPara.Range.Find.Execute("Alex Thurston", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Jay Peterson", 2)

# --- Date: 03/17 -> 4/07 ---------------------------------------------------
# Rebuild the whole paragraph so the stray proofing marks are gone and the
# new date digits each land in their own underlined run.
$datePara = $d.Paragraphs.Item(4)
$dateXml = '<w:p' + $wNs + '>' + `
    '<w:r><w:t>Date:___</w:t></w:r>' + `
    '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>4</w:t></w:r>' + `
    '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>/</w:t></w:r>' + `
    '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>0</w:t></w:r>' + `
    '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>7</w:t></w:r>' + `
    '<w:r><w:t>___</w:t></w:r>' + `
    '</w:p>'
[void]$datePara.Range.InsertXML($dateXml)

# --- Status table updates ---------------------------------------------------
$table = $d.Tables.Item(1)

# Row 2 - Alex Thurston
$table.Cell(2, 2).Range.Text = "Lots of styling changes for each page."
$table.Cell(2, 3).Range.Text = "Keep updating pages and making them look better. Finalizing UI and resource compatibility."
$table.Cell(2, 4).Range.Text = "N/A"

# Row 3 - Nathan Merrill
$table.Cell(3, 2).Range.Text = "Figured out user authentication and account creation. Figured out event creation."
$table.Cell(3, 3).Range.Text = "Supervisor and lot attendant permissions."
$table.Cell(3, 4).Range.Text = "Lots of ways to verify customers. Need to figure out a reliable way."

# Row 4 - Jay Peterson (these cells start out split across multiple runs,
# so rebuild them via InsertXML to collapse everything into one clean run)
[void]$table.Cell(4, 2).Range.InsertXML('<w:p' + $wNs + '><w:r><w:t xml:space="preserve">Figured out user profile editing. </w:t></w:r></w:p>')
[void]$table.Cell(4, 3).Range.InsertXML('<w:p' + $wNs + '><w:r><w:t>Get the password update working. Make html look nice. Unit testing. Verification.</w:t></w:r></w:p>')
$table.Cell(4, 4).Range.Text = "Lots of ways to verify customers. Need to figure out a reliable way."
